# Apply the edit described by the diff:
#  - Insert 2 new data rows at row 224/225 (pushing the existing rows 224-306
#    down to 226-308).
#  - Populate the two new rows with their data.
#
# Net effect on dimension: A1:T306 -> A1:T308

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two blank rows above row 224; everything that was in rows 224-306
# shifts down to rows 226-308.
$ws.Range("224:225").Insert()

# New row 224 data
$ws.Cells.Item(224, 1).Value = 6
$ws.Cells.Item(224, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(224, 3).Value = "Metropolitana"
$ws.Cells.Item(224, 4).Value = 44985
$ws.Cells.Item(224, 5).Value = 13
$ws.Cells.Item(224, 6).Value = "Fruta"
$ws.Cells.Item(224, 7).Value = 100101
$ws.Cells.Item(224, 8).Value = "Berries"
$ws.Cells.Item(224, 9).Value = 100101004
$ws.Cells.Item(224, 10).Value = "Frambuesa"
$ws.Cells.Item(224, 11).Value = "Sin especificar"
$ws.Cells.Item(224, 12).Value = "Primera"
$ws.Cells.Item(224, 13).Value = 250
$ws.Cells.Item(224, 14).Value = 6000
$ws.Cells.Item(224, 15).Value = 6000
$ws.Cells.Item(224, 16).Value = 6000
$ws.Cells.Item(224, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(224, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(224, 19).Value = 3000
$ws.Cells.Item(224, 20).Value = 2

# New row 225 data
$ws.Cells.Item(225, 1).Value = 6
$ws.Cells.Item(225, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(225, 3).Value = "Metropolitana"
$ws.Cells.Item(225, 4).Value = 44985
$ws.Cells.Item(225, 5).Value = 13
$ws.Cells.Item(225, 6).Value = "Fruta"
$ws.Cells.Item(225, 7).Value = 100101
$ws.Cells.Item(225, 8).Value = "Berries"
$ws.Cells.Item(225, 9).Value = 100101004
$ws.Cells.Item(225, 10).Value = "Frambuesa"
$ws.Cells.Item(225, 11).Value = "Sin especificar"
$ws.Cells.Item(225, 12).Value = "Segunda"
$ws.Cells.Item(225, 13).Value = 210
$ws.Cells.Item(225, 14).Value = 4500
$ws.Cells.Item(225, 15).Value = 4500
$ws.Cells.Item(225, 16).Value = 4500
$ws.Cells.Item(225, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(225, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(225, 19).Value = 2250
$ws.Cells.Item(225, 20).Value = 2
